$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the header row with a leading apostrophe so Excel stores it as
# explicit text (quotePrefix) even though the values are unchanged.
$ws.Range("A1").Value = "'UserID"
$ws.Range("B1").Value = "'UserName"
$ws.Range("C1").Value = "'FirstName"
$ws.Range("D1").Value = "'LastName"
$ws.Range("E1").Value = "'Email"
$ws.Range("F1").Value = "'Password"
$ws.Range("G1").Value = "'Phone"

# The UserID column (A2:A4) switches from numeric values to text values
# (quote-prefixed, with extra digits) - logs/ids reformatted as text.
$ws.Range("A2").Value = "'1010345"
$ws.Range("A3").Value = "'10225"
$ws.Range("A4").Value = "'1030885"

# Move the active selection to E9, matching the saved cursor position.
$ws.Range("E9").Select() | Out-Null
